# Applies the cryptocurrency price/volume refresh described by the commit diff.
# (GitHub Actions scheduled update of cryptos.xlsx.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text (e.g. "212.85") must be forced to Text so Excel does not
# silently reinterpret it as a Number; a leading apostrophe is the standard COM/
# Excel "treat the following as text" marker and is not itself stored in the cell.

# Row 2: D2, E2
$ws.Range('D2').Value = '29.686.45'
$ws.Range('E2').Value = '  +0.66%  '

# Row 3: D3, E3
$ws.Range('D3').Value = '1.616.11'
$ws.Range('E3').Value = '  +0.85%  '

# Row 4: E4
$ws.Range('E4').Value = '  -0.50%  '

# Row 5: D5, E5
$ws.Range('D5').Value = '''212.85'
$ws.Range('E5').Value = '  +0.09%  '

# Row 6: E6
$ws.Range('E6').Value = '  -0.41%  '

# Row 7: D7, E7
$ws.Range('D7').Value = '''0.993'
$ws.Range('E7').Value = '  -0.52%  '

# Row 8: D8, E8
$ws.Range('D8').Value = '''28.96'
$ws.Range('E8').Value = '  +7.80%  '

# Row 9: E9
$ws.Range('E9').Value = '  +3.27%  '

# Row 10: D10, E10
$ws.Range('D10').Value = '''0.0609'
$ws.Range('E10').Value = '  +1.76%  '

# Row 11: D11, E11
$ws.Range('D11').Value = '''0.0908'
$ws.Range('E11').Value = '  -0.37%  '

# Row 12: D12, E12
$ws.Range('D12').Value = '1.843.46'
$ws.Range('E12').Value = '  +0.62%  '

# Row 13: D13, E13
$ws.Range('D13').Value = '1.611.06'
$ws.Range('E13').Value = '  +0.72%  '

# Row 14: E14
$ws.Range('E14').Value = '  +6.09%  '

# Row 15: E15
$ws.Range('E15').Value = '  +3.68%  '

# Row 16: D16, E16
$ws.Range('D16').Value = '29.695.71'
$ws.Range('E16').Value = '  +0.63%  '

# Row 17: E17
$ws.Range('E17').Value = '  +16.04%  '

# Row 18: D18, E18
$ws.Range('D18').Value = '''64.55'
$ws.Range('E18').Value = '  +1.79%  '

# Row 19: D19, E19
$ws.Range('D19').Value = '''240.97'
$ws.Range('E19').Value = '  -0.76%  '

# Row 20: D20
$ws.Range('D20').Value = '0.0₃0706'

# Row 21: E21
$ws.Range('E21').Value = '  -0.39%  '

# Row 22: D22, E22
$ws.Range('D22').Value = '''4.10'
$ws.Range('E22').Value = '  +2.91%  '

# Row 23: D23, E23
$ws.Range('D23').Value = '''9.69'
$ws.Range('E23').Value = '  +5.76%  '

# Row 24: D24, E24
$ws.Range('D24').Value = '''2.11'
$ws.Range('E24').Value = '  +0.93%  '

# Row 25: D25, E25
$ws.Range('D25').Value = '''156.65'
$ws.Range('E25').Value = '  +1.38%  '

# Row 26: D26, E26
$ws.Range('D26').Value = '''15.67'
$ws.Range('E26').Value = '  +2.44%  '

# Row 27: E27
$ws.Range('E27').Value = '  +1.24%  '

# Row 28: D28, E28
$ws.Range('D28').Value = '''6.59'
$ws.Range('E28').Value = '  +2.96%  '

# Row 29: E29
$ws.Range('E29').Value = '  -0.43%  '

# Row 30: D30, E30
$ws.Range('D30').Value = '''0.0482'
$ws.Range('E30').Value = '  +1.79%  '

# Row 31: D31, E31
$ws.Range('D31').Value = '''3.30'
$ws.Range('E31').Value = '  +2.52%  '

# Row 32: E32
$ws.Range('E32').Value = '  +1.15%  '

# Row 33: D33, E33
$ws.Range('D33').Value = '''3.19'
$ws.Range('E33').Value = '  +3.02%  '

# Row 34: D34, E34
$ws.Range('D34').Value = '1.438.16'
$ws.Range('E34').Value = '  +1.43%  '

# Row 35: D35, E35
$ws.Range('D35').Value = '''1.61'
$ws.Range('E35').Value = '  +5.97%  '

# Row 36: E36
$ws.Range('E36').Value = '  +2.21%  '

# Row 37: E37
$ws.Range('E37').Value = '  +3.85%  '

# Row 38: E38
$ws.Range('E38').Value = '  -0.92%  '

# Row 39: E39
$ws.Range('E39').Value = '  +3.30%  '

# Row 40: E40
$ws.Range('E40').Value = '  +3.44%  '

# Row 41: D41, E41
$ws.Range('D41').Value = '''0.0506'
$ws.Range('E41').Value = '  +5.13%  '

# Row 42: E42
$ws.Range('E42').Value = '  +0.95%  '

# Row 43: D43, E43
$ws.Range('D43').Value = '''0.825'
$ws.Range('E43').Value = '  +4.07%  '

# Row 44: B44, C44, D44, E44
$ws.Range('B44').Value = 'BitcoinSV'
$ws.Range('C44').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D44').Value = '''54.33'
$ws.Range('E44').Value = '  +2.41%  '

# Row 45: B45, C45, D45, E45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''69.62'
$ws.Range('E45').Value = '  +6.13%  '

# Row 46: E46
$ws.Range('E46').Value = '  -0.48%  '

# Row 47: E47
$ws.Range('E47').Value = '  +20.89%  '

# Row 48: D48, E48
$ws.Range('D48').Value = '''5.44'
$ws.Range('E48').Value = '  +3.02%  '

# Row 49: D49, E49
$ws.Range('D49').Value = '1.753.22'
$ws.Range('E49').Value = '  +0.59%  '

# Row 50: E50
$ws.Range('E50').Value = '  +1.21%  '

# Row 51: E51
$ws.Range('E51').Value = '  -0.84%  '
